# Applies the "chore: update Sheets via scheduled runner" edit:
# refreshed market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the per-server Leve-profit worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 1947.814
$ws.Range("I15").Value = 1947.814
$ws.Range("K15").Value = 5843.442
$ws.Range("M15").Value = -5674.442

# Row 19: Unbreak My Heart
$ws.Range("H19").Value = 812
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -325

# Row 32: Automata for the People
$ws.Range("H32").Value = 1399
$ws.Range("I32").Value = 1999.5
$ws.Range("J32").Value = 798.5
$ws.Range("K32").Value = 1999.5
$ws.Range("L32").Value = 798.5
$ws.Range("M32").Value = -1673.5
$ws.Range("N32").Value = -1450.5

# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 2135.6924
$ws.Range("I38").Value = 350
$ws.Range("J38").Value = 3666.2856
$ws.Range("K38").Value = 1050
$ws.Range("L38").Value = 10998.8568
$ws.Range("M38").Value = -678
$ws.Range("N38").Value = -11742.8568

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 3576.973
$ws.Range("I40").Value = 2692.6155
$ws.Range("J40").Value = 5667.273
$ws.Range("K40").Value = 2692.6155
$ws.Range("L40").Value = 5667.273
$ws.Range("M40").Value = -2517.6155
$ws.Range("N40").Value = -6017.273

# Row 55: A Real Smooth Move
$ws.Range("H55").Value = 1965.6666
$ws.Range("I55").Value = 499.5
$ws.Range("J55").Value = 4898
$ws.Range("K55").Value = 499.5
$ws.Range("L55").Value = 4898
$ws.Range("M55").Value = -285.5
$ws.Range("N55").Value = -5326

# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 1389.7
$ws.Range("I58").Value = 79.40000000000001
$ws.Range("K58").Value = 238.2
$ws.Range("M58").Value = -88.20000000000002

# Row 95: Official Strategy Guide
$ws.Range("H95").Value = 34499.5
$ws.Range("J95").Value = 34499.5
$ws.Range("L95").Value = 34499.5
$ws.Range("N95").Value = -39991.5

# Row 96: Scroll Down
$ws.Range("H96").Value = 237.5
$ws.Range("J96").Value = 204.5
$ws.Range("L96").Value = 613.5
$ws.Range("N96").Value = -3359.5

# Row 97: Materia Worth
$ws.Range("H97").Value = 995
$ws.Range("J97").Value = 995
$ws.Range("L97").Value = 2985
$ws.Range("N97").Value = -3977

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 2201.9167
$ws.Range("I100").Value = 1945.9231
$ws.Range("J100").Value = 2504.4546
$ws.Range("K100").Value = 1945.9231
$ws.Range("L100").Value = 2504.4546
$ws.Range("M100").Value = -1404.9231
$ws.Range("N100").Value = -3586.4546

# Row 113: Amaro Kart
$ws.Range("H113").Value = 4999.5

# Row 125: Body over Mind
$ws.Range("H125").Value = 3499.375
$ws.Range("I125").Value = 2333
$ws.Range("J125").Value = 4199.2
$ws.Range("K125").Value = 20997
$ws.Range("L125").Value = 37792.8
$ws.Range("M125").Value = -18537
$ws.Range("N125").Value = -42712.8

# Row 131: Mindful Study
$ws.Range("H131").Value = 3895.7144
$ws.Range("I131").Value = 4054
$ws.Range("J131").Value = 3500
$ws.Range("K131").Value = 12162
$ws.Range("L131").Value = 10500
$ws.Range("M131").Value = -7122
$ws.Range("N131").Value = -20580

# Row 135: For Tired Minds
$ws.Range("H135").Value = 1010.4286
$ws.Range("I135").Value = 758.6667
$ws.Range("K135").Value = 6828.0003
$ws.Range("M135").Value = -4293.0003

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1929.75
$ws.Range("I137").Value = 1033.4615
$ws.Range("K137").Value = 3100.3845
$ws.Range("M137").Value = -550.3844999999997

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 145.21428
$ws.Range("I5").Value = 104.125
$ws.Range("K5").Value = 104.125
$ws.Range("M5").Value = 7.875

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 5915.857
$ws.Range("I74").Value = 6066.1816
$ws.Range("J74").Value = 5364.6665
$ws.Range("K74").Value = 6066.1816
$ws.Range("L74").Value = 5364.6665
$ws.Range("M74").Value = -5192.1816
$ws.Range("N74").Value = -7112.6665

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 5915.857
$ws.Range("I77").Value = 6066.1816
$ws.Range("J77").Value = 5364.6665
$ws.Range("K77").Value = 30330.908
$ws.Range("L77").Value = 26823.3325
$ws.Range("M77").Value = -25962.908
$ws.Range("N77").Value = -35559.3325

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 2354.9285
$ws.Range("I110").Value = 1755.5
$ws.Range("K110").Value = 1755.5
$ws.Range("M110").Value = 289.5

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2694.75
$ws.Range("I122").Value = 2732.3333
$ws.Range("K122").Value = 8196.999899999999
$ws.Range("M122").Value = -5746.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Range("H4").Value = 145.21428
$ws.Range("I4").Value = 104.125
$ws.Range("K4").Value = 104.125
$ws.Range("M4").Value = 10.875

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 3889.3333
$ws.Range("I86").Value = 1603.4445
$ws.Range("J86").Value = 7318.1665
$ws.Range("K86").Value = 1603.4445
$ws.Range("L86").Value = 7318.1665
$ws.Range("M86").Value = -480.4445000000001
$ws.Range("N86").Value = -9564.166499999999

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 3889.3333
$ws.Range("I89").Value = 1603.4445
$ws.Range("J89").Value = 7318.1665
$ws.Range("K89").Value = 8017.2225
$ws.Range("L89").Value = 36590.8325
$ws.Range("M89").Value = -2401.2225
$ws.Range("N89").Value = -47822.8325

# Row 94: High Steal
$ws.Range("H94").Value = 562.6667
$ws.Range("I94").Value = 494.16666
$ws.Range("K94").Value = 494.16666
$ws.Range("M94").Value = -43.16665999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 2963.1667
$ws.Range("J22").Value = 3887.5
$ws.Range("L22").Value = 3887.5
$ws.Range("N22").Value = -4587.5

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 4046.4285
$ws.Range("I132").Value = 4046.4285
$ws.Range("K132").Value = 12139.2855
$ws.Range("M132").Value = -9609.2855

$ws = $wb.Worksheets.Item("CUL")
# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 346.625
$ws.Range("I38").Value = 357.9524
$ws.Range("J38").Value = 267.33334
$ws.Range("K38").Value = 1073.8572
$ws.Range("L38").Value = 802.0000200000001
$ws.Range("M38").Value = -726.8571999999999
$ws.Range("N38").Value = -1496.00002

# Row 86: Let's Not Get Sappy
$ws.Range("H86").Value = 449.75
$ws.Range("I86").Value = 120
$ws.Range("J86").Value = 779.5
$ws.Range("K86").Value = 360
$ws.Range("L86").Value = 2338.5
$ws.Range("M86").Value = 826
$ws.Range("N86").Value = -4710.5

# Row 89: Luxury Spillover (L)
$ws.Range("H89").Value = 449.75
$ws.Range("I89").Value = 120
$ws.Range("J89").Value = 779.5
$ws.Range("K89").Value = 1080
$ws.Range("L89").Value = 7015.5
$ws.Range("M89").Value = 4848
$ws.Range("N89").Value = -18871.5

# Row 107: Slippery Service
$ws.Range("H107").Value = 477.8421
$ws.Range("I107").Value = 404.2
$ws.Range("K107").Value = 1212.6
$ws.Range("M107").Value = 707.4000000000001

# Row 109: Cure for What Ails
$ws.Range("H109").Value = 842.4545000000001
$ws.Range("I109").Value = 846.3333
$ws.Range("K109").Value = 2538.9999
$ws.Range("M109").Value = -1498.9999

# Row 132: More Mezcal
$ws.Range("H132").Value = 1150
$ws.Range("I132").Value = 947
$ws.Range("K132").Value = 8523
$ws.Range("M132").Value = -5993

$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success
$ws.Range("H11").Value = 14696217
$ws.Range("I11").Value = 12353706
$ws.Range("K11").Value = 12353706
$ws.Range("M11").Value = -12353567

# Row 57: Gold Is So Last Year
$ws.Range("H57").Value = 26248.75
$ws.Range("I57").Value = 8331.666999999999
$ws.Range("J57").Value = 80000
$ws.Range("K57").Value = 8331.666999999999
$ws.Range("L57").Value = 80000
$ws.Range("M57").Value = -7511.666999999999
$ws.Range("N57").Value = -81640

# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 199.75
$ws.Range("I107").Value = 194.4
$ws.Range("J107").Value = 208.66667
$ws.Range("K107").Value = 194.4
$ws.Range("L107").Value = 208.66667
$ws.Range("M107").Value = 1725.6
$ws.Range("N107").Value = -4048.66667

# Row 135: Fan of the Foreign
$ws.Range("H135").Value = 225390
$ws.Range("J135").Value = 225390
$ws.Range("L135").Value = 225390
$ws.Range("N135").Value = -235530

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 8001
$ws.Range("I16").Value = 8001
$ws.Range("K16").Value = 8001
$ws.Range("M16").Value = -7831

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 5559.7144
$ws.Range("I46").Value = 2229.5
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 2229.5
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -2041.5
$ws.Range("N46").Value = -10376

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1274.091
$ws.Range("I93").Value = 1242.5883
$ws.Range("K93").Value = 1242.5883
$ws.Range("M93").Value = 5.41170000000011

$ws = $wb.Worksheets.Item("WVR")
# Row 52: Party Animals
$ws.Range("H52").Value = 18985.25
$ws.Range("J52").Value = 45999
$ws.Range("L52").Value = 45999
$ws.Range("N52").Value = -46451

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 3689.625
$ws.Range("I122").Value = 1379.25
$ws.Range("K122").Value = 4137.75
$ws.Range("M122").Value = -1687.75
